$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "Asistido Copilot" (column D) values for Bloque 1 (rows 6-9)
$ws.Range("D6").Value = 3
$ws.Range("D7").Value = 3
$ws.Range("D8").Value = 4
$ws.Range("D9").Value = 1

# Fill in "Asistido Copilot" (column D) values for Bloque 2 (rows 13-16)
$ws.Range("D13").Value = 4
$ws.Range("D14").Value = 4
$ws.Range("D15").Value = 7
$ws.Range("D16").Value = 6

# Update the view to match the saved selection/scroll position
$ws.Range("D17").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
